$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml($matchPrefix, $innerXml) {
    $target = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($matchPrefix)) {
            $target = $p
            break
        }
    }
    if ($null -eq $target) {
        throw "Paragraph starting with '$matchPrefix' not found"
    }
    $xml = '<w:p ' + $wNs + '>' + $innerXml + '</w:p>'
    [void]$target.Range.InsertXML($xml)
}

# Paragraph 1: Portuguese "Programa" listing
Set-ParagraphXml "1) Medidores. Osciloscópio." (
    '<w:r>' +
        '<w:t>1) Medidores. Osciloscópio.</w:t><w:br/>' +
        '<w:t>2) Tensão alternada.</w:t><w:br/>' +
        '<w:t xml:space="preserve">3) Potências. </w:t><w:br/>' +
        '<w:t>4) Filtros.</w:t><w:br/>' +
        '<w:t>5) Ressonância.</w:t><w:br/>' +
        '<w:t>6) Campo magnético alternado.</w:t>' +
    '</w:r>'
)

# Paragraph 2: English "Programa" listing (italic)
Set-ParagraphXml "1) Meters. Oscilloscope." (
    '<w:r>' +
        '<w:rPr><w:i/></w:rPr>' +
        '<w:t>1) Meters. Oscilloscope.</w:t><w:br/>' +
        '<w:t>2) AC voltage.</w:t><w:br/>' +
        '<w:t>3) Powers.</w:t><w:br/>' +
        '<w:t>4) Filters.</w:t><w:br/>' +
        '<w:t>5) Resonance.</w:t><w:br/>' +
        '<w:t>6) AC Magnetic fields.</w:t>' +
    '</w:r>'
)

# Paragraph 3: Bibliografia
Set-ParagraphXml "CAPUANO, G. Francisco" (
    '<w:r>' +
        '<w:t xml:space="preserve">CAPUANO, G. Francisco; MARINO, M.A. Maria. Laboratório de eletricidade </w:t><w:br/>' +
        '<w:t>Eletrônica, Editora Érica (1998).</w:t><w:br/>' +
        '<w:t>MARKUS, Otávio. Circuitos elétricos: corrente contínua e corrente alternada-</w:t><w:br/>' +
        '<w:t>Teoria e Exercícios, Editora Érica, (2008).</w:t><w:br/>' +
        '<w:t>SADIKU, Mathew N. O.; ALEXANDER, Charles. Fundamentos de circuitos elétricos, Mcgraw-hill Interamericana (2009).</w:t>' +
    '</w:r>'
)

Write-Output "Done"
